$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 7.915999999999999
$ws.Range("D5").Value = -8.253
$ws.Range("A10").Value = -21.653
$ws.Range("A12").Value = -21.303
$ws.Range("B12").Value = 6.303999999999999
$ws.Range("C12").Value = -11.25
$ws.Range("C13").Value = -12.857
$ws.Range("B17").Value = 5.055
$ws.Range("A18").Value = -21.916
$ws.Range("C21").Value = -12.907
$ws.Range("D23").Value = -8.356999999999999
$ws.Range("B26").Value = 5.501
$ws.Range("D26").Value = -7.798
$ws.Range("B27").Value = 5.789000000000001
$ws.Range("B28").Value = 5.390000000000001
$ws.Range("D34").Value = -7.739999999999999
$ws.Range("C36").Value = -12.776
$ws.Range("A37").Value = -21.206
$ws.Range("B37").Value = 7.284999999999999
$ws.Range("C38").Value = -12.4
$ws.Range("D40").Value = -8.125
$ws.Range("C41").Value = -12.554
$ws.Range("D47").Value = -7.699
$ws.Range("C52").Value = -11.775
$ws.Range("A55").Value = -22.109
$ws.Range("D60").Value = -8.198
$ws.Range("B65").Value = 5.823
$ws.Range("C67").Value = -10.807
$ws.Range("A68").Value = -21.502
$ws.Range("D72").Value = -7.478
$ws.Range("B73").Value = 6.543000000000001
$ws.Range("A77").Value = -21.11
$ws.Range("A78").Value = -20.752
$ws.Range("D83").Value = -7.834000000000001
$ws.Range("B84").Value = 5.659000000000001
$ws.Range("B85").Value = 4.975
$ws.Range("C89").Value = -13.634
$ws.Range("B93").Value = 5.587000000000001
$ws.Range("B95").Value = 6.603
$ws.Range("C95").Value = -12.118
$ws.Range("B98").Value = 7.281000000000001
$ws.Range("B99").Value = 5.231
$ws.Range("B101").Value = 6.360999999999999
$ws.Range("C105").Value = -12.753
